# Correct the condition/count for the concatenate script and the
# Uncertainty Table entry where X_CAT = Average Distance.
# The row for parameter L_E_12 / "Average Distance" on the Unified_table
# sheet was mis-numbered as 22; it should be 23. All of the subsequent
# rows recompute their running index via formulas that chain off this
# cell (=A43+1, etc.), so updating this one value cascades through the
# rest of the numbering column automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unified_table")

$ws.Range("A43").Value = 23

$excel.CalculateFullRebuild()
